$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-18 08:23:01"
$wsZh.Range("G2").Value = "2016-02-18 08:23:44"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-18 08:23:15"
$wsDe.Range("G2").Value = "2016-02-18 08:24:07"
